$wb = $excel.ActiveWorkbook

# --- Rename the third sheet (Sheet1 -> doSaveUpdates) and populate it ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "doSaveUpdates"

# Header row
$ws.Range("A1").Value = "phoneNumber"
$ws.Range("B1").Value = "streetAddress"
$ws.Range("C1").Value = "city"

# Fill column B (street) and C (city) in the order that reproduces the
# original shared-string table layout, then column A (phone numbers).
$ws.Range("B4").Value = "125 Main st"
$ws.Range("B2").Value = "123 Main st"
$ws.Range("B3").Value = "124 Main st"

$ws.Range("C2").Value = "Philadelphia"
$ws.Range("C3").Value = "New York"
$ws.Range("C4").Value = "Pitsburgh"

$ws.Range("A2").Value = 2679875852
$ws.Range("A3").Value = 2678526547
$ws.Range("A4").Value = 2674718956

# Column widths (closest attainable values to the authored widths)
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 13.33203125

# Select D8 on this sheet (becomes the sheet's stored selection / cursor)
$ws.Range("D8").Select() | Out-Null

# --- Move the active tab back to "signin" ---
$wsSignin = $wb.Worksheets.Item("signin")
$wsSignin.Activate() | Out-Null
